$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1 (A1) with today's rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.78 = 35575.15 pesos`n✅ 35575.15 pesos = 8.74 = 910.02 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate table on the "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 113.9
$ws2.Range("O10").Value = 4052.01
$ws2.Range("N12").Value = 4069.99
$ws2.Range("O12").Value = 104.111
